$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.823.09"
$ws.Range("E2").Value = "  -0.88%  "
$ws.Range("D3").Value = "2.033.79"
$ws.Range("E3").Value = "  -1.22%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'227.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.10%  "
$ws.Range("E6").Value = "  -0.84%  "
$ws.Range("D7").Value = "'60.44"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.80%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "'0.378"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.28%  "
$ws.Range("E10").Value = "  +0.58%  "
$ws.Range("D11").Value = "'0.103"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.30%  "
$ws.Range("D12").Value = "2.336.47"
$ws.Range("E12").Value = "  -0.97%  "
$ws.Range("D13").Value = "'14.48"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.90%  "
$ws.Range("D14").Value = "'21.31"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.13%  "
$ws.Range("D15").Value = "'0.761"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.81%  "
$ws.Range("E16").Value = "  -2.33%  "
$ws.Range("D17").Value = "2.051.52"
$ws.Range("E17").Value = "  -0.66%  "
$ws.Range("D18").Value = "37.788.05"
$ws.Range("E18").Value = "  -0.64%  "
$ws.Range("D19").Value = "'69.80"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.14%  "
$ws.Range("D20").Value = "'5.89"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.80%  "
$ws.Range("E21").Value = "  -1.71%  "
$ws.Range("D22").Value = "'223.94"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.40%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("D26").Value = "'9.39"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.89%  "
$ws.Range("D27").Value = "'167.33"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.46%  "
$ws.Range("D28").Value = "'0.129"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.60%  "
$ws.Range("D29").Value = "'18.87"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.96%  "
$ws.Range("E30").Value = "  -3.92%  "
$ws.Range("E31").Value = "  +0.38%  "
$ws.Range("D32").Value = "'2.25"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +9.01%  "
$ws.Range("D33").Value = "'4.40"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.42%  "
$ws.Range("D34").Value = "'0.0606"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.23%  "
$ws.Range("D35").Value = "'4.50"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.19%  "
$ws.Range("D36").Value = "'6.39"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.63%  "
$ws.Range("D37").Value = "'2.29"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.71%  "
$ws.Range("D38").Value = "'3.34"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.15%  "
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("D40").Value = "'17.67"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.81%  "
$ws.Range("D41").Value = "1.537.80"
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("E42").Value = "  +0.19%  "
$ws.Range("D43").Value = "'96.31"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.09%  "
$ws.Range("E44").Value = "  -2.46%  "
$ws.Range("D45").Value = "'0.0913"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.26%  "
$ws.Range("E46").Value = "  -3.12%  "
$ws.Range("D47").Value = "'3.99"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.00%  "
$ws.Range("E48").Value = "  -1.36%  "
$ws.Range("D49").Value = "'2.96"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.26%  "
$ws.Range("D50").Value = "'7.12"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.38%  "
$ws.Range("D51").Value = "2.226.15"
$ws.Range("E51").Value = "  -0.90%  "
